# Updated cryptos list on Fri Mar 24 08:50:58 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) for each coin row, and swaps the
# WEMIXTOKEN / FraxShare rows (42 <-> 43) per the latest ranking pull.
#
# Note: Price values that look like plain decimals (e.g. "1.006") are
# written with a leading apostrophe so Excel keeps them as text (matching
# the original inlineStr cell type) instead of silently converting them to
# numbers. Values that already contain two dots (e.g. "28.080.99") are
# unambiguous text and don't need the apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.080.99"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.802.89"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'324.44"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4278"
$ws.Range("E7").Value = "  -3.97%  "
$ws.Range("D8").Value = "'0.3626"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").Value = "'44.75"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'0.07588"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'1.128"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("D13").Value = "'21.61"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "'6.235"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "'7.389"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "1.823.52"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("D17").Value = "'93.10"
$ws.Range("E17").Value = "  +5.76%  "
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'0.06366"
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'17.26"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").Value = "'6.176"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "28.137.54"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "'11.48"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'2.140"
$ws.Range("E25").Value = "  -7.55%  "
$ws.Range("D26").Value = "'160.39"
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").Value = "'20.52"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "2.029.31"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("D29").Value = "'2.231"
$ws.Range("E29").Value = "  -5.87%  "
$ws.Range("D30").Value = "'129.36"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'1.181"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").Value = "'5.908"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").Value = "'0.09047"
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").Value = "'3.533"
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("D35").Value = "'12.82"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "'0.02386"
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").Value = "'5.141"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "'0.6531"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").Value = "'0.2128"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("D40").Value = "'0.06128"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").Value = "'1.201"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.979"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.424"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'13.72"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "'0.6015"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'3.729"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "'125.16"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").Value = "'1.998"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "'1.163"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").Value = "'0.06981"
$ws.Range("E51").Value = "  +1.06%  "
